# [UPDATE] Datos conexion BD + Renombre objetos Alta de Cuentas
#
# Sheet "Database" holds Name/Value pairs:
#   A4 = dBUser  / B4 = <username>
#   A5 = dBPass  / B5 = <password>
#
# Update dBUser and dBPass to new rotated credentials, then leave the
# selection on B4 (the last-edited cell) instead of B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Database")

$ws.Range("B4").Value = "kqg4sm1ik53xfqvs984d"
$ws.Range("B5").Value = "pscale_pw_8lmwdeLbOlpwgLrIVyJNwhmmNZfq3xTokRdG8IuZ34O"

$ws.Range("B4").Select()
